$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the old "Total:" row (row 29) so it shifts down to row 30,
# making room for one more actual/task entry (rows 25-29 now hold 5 entries).
$ws.Rows(29).Insert()

# Row 25: task relabeled to "19-David"; its point now logged under Points (M) instead of Days (L)
$ws.Range("J25").Value = "19-David"
$ws.Range("L25").Clear()
$ws.Range("M25").Value = 1

# Row 27: old "24a" task (now done by Dann) moves here with its Days value
$ws.Range("J27").Value = "24a-Dann"
$ws.Range("O27").Clear()
$ws.Range("L27").Value = 3

# Row 28: old "24b" task (now done by Dann) moves here with its Days value
$ws.Range("J28").Value = "24b-Dann"
$ws.Range("O28").Clear()
$ws.Range("L28").Value = 3

# Row 26: new task "22-Riaz" logged under Total Story Points (O)
$ws.Range("J26").Value = "22-Riaz"
$ws.Range("L26").Clear()
$ws.Range("O26").Value = 5

# Row 29 (newly inserted, previously blank): new task "33-Philip"
$ws.Range("J29").Value = "33-Philip"
$ws.Range("O29").Value = 8

# Row 30 (old Total row, shifted down by the insert): extend the sum range
$ws.Range("K30").Formula = "=SUM(K25:O29)"

# Row 15 actual-burn-down formulas: extend each SUM range to include the new row 29
$ws.Range("D15").Formula = "=C15-SUM(K25:K29)"
$ws.Range("E15").Formula = "=D15-SUM(L25:L29)"
$ws.Range("F15").Formula = "=E15-SUM(M25:M29)"
$ws.Range("G15").Formula = "=F15-SUM(N25:N29)"
$ws.Range("H15").Formula = "=G15-SUM(O25:O29)"

# Restore the view: scrolled down a bit, with K30 (the new total cell) selected
$excel.ActiveWindow.ScrollRow = 4
$ws.Range("K30").Select() | Out-Null

$wb.Application.Calculate() | Out-Null
